$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 11
$ws.Range("F3").Value  = 11441
$ws.Range("F4").Value  = 1305
$ws.Range("F8").Value  = 175
$ws.Range("F9").Value  = 965
$ws.Range("F11").Value = 2335
$ws.Range("F13").Value = 1156
$ws.Range("F16").Value = 869
$ws.Range("F17").Value = 1042
$ws.Range("F21").Value = 728
$ws.Range("F25").Value = 71
$ws.Range("F26").Value = 496
$ws.Range("F27").Value = 551
$ws.Range("F29").Value = 276
$ws.Range("F31").Value = 647
$ws.Range("F32").Value = 2772
$ws.Range("F33").Value = 444
$ws.Range("F34").Value = 39
$ws.Range("F35").Value = 293
$ws.Range("F37").Value = 98
$ws.Range("F38").Value = 1534

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 135

# Sheet "本地生活" (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2238
$ws.Range("F3").Value = 691

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 2238
$ws.Range("F3").Value  = 691
$ws.Range("F4").Value  = 11
$ws.Range("F5").Value  = 11441
$ws.Range("F10").Value = 965
$ws.Range("F11").Value = 2335
$ws.Range("F13").Value = 1156
$ws.Range("F16").Value = 869
$ws.Range("F17").Value = 1042
$ws.Range("F23").Value = 728
$ws.Range("F28").Value = 71
$ws.Range("F29").Value = 496
$ws.Range("F30").Value = 551
$ws.Range("F33").Value = 2773
$ws.Range("F35").Value = 444
$ws.Range("F36").Value = 39
$ws.Range("F38").Value = 98
$ws.Range("F39").Value = 1534
